# "updated map and added Collector command"
# Adds a new part row (Ball Collector Motor, port 5, comment about the talon)
# to the map worksheet, immediately below the existing USB webcam row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Ball Collector Motor"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "Needs to be a talon hooked up to port 5 on the rio"

# Move the active selection onto the newly added row, matching the
# author's cursor position after entering the data.
$ws.Range("C7").Select()
